# Generate Report for Handback
# Applies the "handback" report values to the localization-status workbook:
#  - Updates the "In Translation" status text to "Handed back: in sync with en-US"
#    everywhere it appears (Overview Priority/Content-Duplicate columns as well as
#    the Status column on the language sheets).
#  - Fills in the previously-empty "Latest Target File" / "Latest Handback File"
#    columns (and, for de-de, the "Latest Handback DateTime") now that a handback
#    report has been generated, turning the target-file cell into a hyperlink.
#  - Updates the stale placeholder handback datetime to the real handback time.
#  - Widens a few columns so the newly-populated long file names are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: refresh status text + widen the zh-cn / de-de status columns
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = "8c04d975-2573-41b0-802d-8fdc6235037e.md"
$zhcn.Range("J2").Value = "8c04d975-2573-41b0-802d-8fdc6235037e.6c12181e7d46d3c1e46a585e251d8ad36fcc0160.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 16:24:33"

$zhcn.Range("I3").Value = "9f7221d1-1144-44fa-8aae-a6a57991399c.md"
$zhcn.Range("J3").Value = "9f7221d1-1144-44fa-8aae-a6a57991399c.b9650f7a644d61ab76cfbc31907036676b804a49.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-24 16:24:33"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8e86f8d352c3e59e8a3bcaf7df1f740e752a5c5/e2e/8c04d975-2573-41b0-802d-8fdc6235037e.md", "", "", "8c04d975-2573-41b0-802d-8fdc6235037e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8e86f8d352c3e59e8a3bcaf7df1f740e752a5c5/e2e/9f7221d1-1144-44fa-8aae-a6a57991399c.md", "", "", "9f7221d1-1144-44fa-8aae-a6a57991399c.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = "8c04d975-2573-41b0-802d-8fdc6235037e.md"
$dede.Range("J2").Value = "8c04d975-2573-41b0-802d-8fdc6235037e.6c12181e7d46d3c1e46a585e251d8ad36fcc0160.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 16:24:40"

$dede.Range("I3").Value = "9f7221d1-1144-44fa-8aae-a6a57991399c.md"
$dede.Range("J3").Value = "9f7221d1-1144-44fa-8aae-a6a57991399c.b9650f7a644d61ab76cfbc31907036676b804a49.de-de.xlf"
$dede.Range("K3").Value = "2016-08-24 16:24:40"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8e86f8d352c3e59e8a3bcaf7df1f740e752a5c5/e2e/8c04d975-2573-41b0-802d-8fdc6235037e.md", "", "", "8c04d975-2573-41b0-802d-8fdc6235037e.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8e86f8d352c3e59e8a3bcaf7df1f740e752a5c5/e2e/9f7221d1-1144-44fa-8aae-a6a57991399c.md", "", "", "9f7221d1-1144-44fa-8aae-a6a57991399c.md")

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
